# Applies the edits described by the commit diff to the active document.
$d = $word.ActiveDocument

# --- 1. Simple text replacements, scoped to the specific paragraph so that
#        duplicate/overlapping text elsewhere in the document is never
#        touched by accident. Paragraph numbers below refer to the ORIGINAL
#        (pre-edit) paragraph layout; none of the edits in this section
#        insert or delete paragraphs, so the numbering stays stable
#        throughout this block. ---

function Replace-InParagraph($index, $find, $replace) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-InParagraph 8 `
    "The website should accept a keyword query and an advanced search that includes all NFL attributes including: height, weight, combine stats, position, wingspan, age and team  " `
    "The website should accept a keyword query and an advanced search that includes all NBA attributes including: name, height, weight, stats, position, headshots (if applicable), "

Replace-InParagraph 10 `
    "The site will (by default) return all players that match the advanced search criteria in order from most relevant to least" `
    "The site will include a full roster with all player stat averages along with the ability to filter/sort through each attribute/stat"

Replace-InParagraph 12 `
    "The site will allow users to sort players by stats and/or attributes (through row/column roster layout) " `
    "Site will have a “Highlights” page to display current player news and weekly highlights"

Replace-InParagraph 14 `
    "Site will have a “Seasonal Leaders” page to highlight the most dominant players" `
    "Website will include full season schedule with broadcast links and/or information "

Replace-InParagraph 16 `
    "Clicking on a player will open up a “Player Information” page that shows seasonal highlights, biography, college information, and full description (attributes and statistics) of selected player " `
    "User Interface will be easily navigable and responsive to allow for a smooth browsing experience "

Replace-InParagraph 18 `
    "Website will include full season schedule with broadcast links and/or information " `
    "User accounts/registration - Users should be able to enter credentials, the form data will be validated and then stored inside the mongodb database"

Replace-InParagraph 25 `
    "NFL Roster API for player stats and attributes. Site will have its on dedicated API that automatically pulls from databases on a weekly basis" `
    "NBA Roster API for player stats, pictures, news/highlights, and attributes."

Replace-InParagraph 27 `
    "Tools and Technology: Javascript, HTML, CSS, MongoDb, Express.js, Angular.js, Node.js, Heroku, Bootstrap (MEAN Stack Web Development) Visual Studio Code, Codepen" `
    "Tools and Technology: Javascript, HTML, CSS, MongoDb, Express.js, React.js, Node.js, Heroku, Bootstrap (MEAN Stack Web Development) Visual Studio Code"

Replace-InParagraph 31 `
    "Website to be fully operational by 8/2" `
    "Website to be fully operational by 7/26"

# --- 2. Indentation tweaks on the two blank spacer paragraphs that sit
#        between the "Highlights" bullet / the "User Interface" bullet
#        and their following bullet item. ---

$d.Paragraphs.Item(13).Format.LeftIndent = 72
$d.Paragraphs.Item(17).Format.LeftIndent = 36

# --- 3. Remove the two trailing paragraphs that are no longer present
#        (a blank spacer and the old duplicate "User Interface..." bullet).
#        Paragraph 19 is deleted twice in a row because after the first
#        delete, the following paragraph slides into slot 19. ---

$d.Paragraphs.Item(19).Range.Delete() | Out-Null
$d.Paragraphs.Item(19).Range.Delete() | Out-Null
